$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2017-02-09 16:10:18"

$wsZhCn.Range("H2").Value = "2017-02-09 16:09:53"
$wsZhCn.Range("L2").Value = "2017-02-09 16:10:57"

$wsDeDe.Range("H2").Value = "2017-02-09 16:10:18"
$wsDeDe.Range("L2").Value = "2017-02-09 16:11:23"
